{"js": "// The \"COMPETENCES TECHNIQUES\" section lists 7 skill lines. The commit\n// re-orders those lines (content moved from a Python list into a JSON\n// structure whose key order differs), while every other paragraph in the\n// document keeps its original text/formatting. We locate the 7\n// consecutive paragraphs by their current text and rewrite each one's\n// text in place to match the new order, so paragraph formatting\n// (w:pPr/spacing etc.) is left completely untouched.\n\nconst oldOrder = [\n  \"Web : api\",\n  \"Langages : r, python, matlab, c, c++\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n  \"Autres : kubernetes, cloud, la conception, github\",\n  \"Visualisation : tableau\",\n  \"ML/AI : fastapi, django, Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n];\n\nconst newOrder = [\n  \"Langages : r, python, matlab, c, c++\",\n  \"Visualisation : tableau\",\n  \"MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Autres : kubernetes, cloud, la conception, github\",\n  \"Web : api\",\n  \"ML/AI : fastapi, django, Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the contiguous block of paragraphs whose text matches oldOrder.\nlet startIndex = -1;\nfor (let i = 0; i + oldOrder.length <= paragraphs.items.length; i++) {\n  let isMatch = true;\n  for (let j = 0; j < oldOrder.length; j++) {\n    if (paragraphs.items[i + j].text !== oldOrder[j]) {\n      isMatch = false;\n      break;\n    }\n  }\n  if (isMatch) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not find the expected 'COMPETENCES TECHNIQUES' skill paragraphs.\");\n}\n\nfor (let j = 0; j < newOrder.length; j++) {\n  paragraphs.items[startIndex + j].insertText(newOrder[j], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The \"COMPETENCES TECHNIQUES\" section lists 7 skill lines. The commit\n# re-orders those lines (content moved from a Python list into a JSON\n# structure whose key order differs), while every other paragraph in the\n# document keeps its original text/formatting. We locate the 7\n# consecutive paragraphs by their current text and rewrite each one's\n# Range.Text in place to match the new order. Assigning Range.Text (as\n# opposed to rebuilding/moving paragraphs) leaves the paragraph mark and\n# w:pPr formatting (spacing, etc.) completely untouched.\n\n$d = $word.ActiveDocument\n\n$oldOrder = @(\n    \"Web : api\",\n    \"Langages : r, python, matlab, c, c++\",\n    \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n    \"Autres : kubernetes, cloud, la conception, github\",\n    \"Visualisation : tableau\",\n    \"ML/AI : fastapi, django, Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n    \"MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n)\n\n$newOrder = @(\n    \"Langages : r, python, matlab, c, c++\",\n    \"Visualisation : tableau\",\n    \"MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n    \"Autres : kubernetes, cloud, la conception, github\",\n    \"Web : api\",\n    \"ML/AI : fastapi, django, Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n    \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n)\n\n$count = $d.Paragraphs.Count\n$startIndex = -1\n\nfor ($i = 1; $i -le ($count - $oldOrder.Length + 1); $i++) {\n    $isMatch = $true\n    for ($j = 0; $j -lt $oldOrder.Length; $j++) {\n        $text = $d.Paragraphs.Item($i + $j).Range.Text\n        # Paragraph ranges include the trailing paragraph mark; strip it\n        # (and any stray cell-mark char) before comparing.\n        $text = $text.TrimEnd([char]13, [char]7)\n        if ($text -ne $oldOrder[$j]) {\n            $isMatch = $false\n            break\n        }\n    }\n    if ($isMatch) {\n        $startIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1) {\n    throw \"Could not locate the expected 'COMPETENCES TECHNIQUES' skill paragraphs.\"\n}\n\nfor ($j = 0; $j -lt $newOrder.Length; $j++) {\n    $d.Paragraphs.Item($startIndex + $j).Range.Text = $newOrder[$j]\n}\n"}
